$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 20249.928
$ws.Range("I51").Value = 4500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4016

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1836.0834
$ws.Range("I86").Value = 1771.7142
$ws.Range("J86").Value = 1926.2
$ws.Range("K86").Value = 1771.7142
$ws.Range("L86").Value = 1926.2
$ws.Range("M86").Value = -648.7141999999999
$ws.Range("N86").Value = -4172.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1836.0834
$ws.Range("I89").Value = 1771.7142
$ws.Range("J89").Value = 1926.2
$ws.Range("K89").Value = 8858.571
$ws.Range("L89").Value = 9631
$ws.Range("M89").Value = -3242.571
$ws.Range("N89").Value = -20863

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1799.5
$ws.Range("J101").Value = 1799.5
$ws.Range("L101").Value = 5398.5
$ws.Range("N101").Value = -8642.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4524.9375
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4524.9375
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 13574.8125
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -15790.8125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5655.75
$ws.Range("I113").Value = 1749.3334
$ws.Range("K113").Value = 1749.3334
$ws.Range("M113").Value = 1504.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2657.1462
$ws.Range("I132").Value = 2756.2368
$ws.Range("J132").Value = 1402
$ws.Range("K132").Value = 8268.7104
$ws.Range("L132").Value = 4206
$ws.Range("M132").Value = -5738.7104
$ws.Range("N132").Value = -9266

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2117.2173
$ws.Range("I137").Value = 1663.8
$ws.Range("J137").Value = 2466
$ws.Range("K137").Value = 4991.4
$ws.Range("L137").Value = 7398
$ws.Range("M137").Value = -2441.4
$ws.Range("N137").Value = -12498

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3986.2163
$ws.Range("I138").Value = 1082.56
$ws.Range("J138").Value = 10035.5
$ws.Range("K138").Value = 3247.68
$ws.Range("L138").Value = 30106.5
$ws.Range("M138").Value = 1892.32
$ws.Range("N138").Value = -40386.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13717.223
$ws.Range("I32").Value = 1149.2609
$ws.Range("J32").Value = 35952.848
$ws.Range("K32").Value = 1149.2609
$ws.Range("L32").Value = 35952.848
$ws.Range("M32").Value = -862.2609
$ws.Range("N32").Value = -36526.848

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3241.75
$ws.Range("I45").Value = 3241.75
$ws.Range("K45").Value = 3241.75
$ws.Range("M45").Value = -2864.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6193.974
$ws.Range("I61").Value = 4043.0312
$ws.Range("K61").Value = 4043.0312
$ws.Range("M61").Value = -3831.0312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2807.2104
$ws.Range("I74").Value = 2521.1875
$ws.Range("K74").Value = 2521.1875
$ws.Range("M74").Value = -1647.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2807.2104
$ws.Range("I77").Value = 2521.1875
$ws.Range("K77").Value = 12605.9375
$ws.Range("M77").Value = -8237.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1349.1578
$ws.Range("J88").Value = 1425.5
$ws.Range("L88").Value = 1425.5
$ws.Range("N88").Value = -2237.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1349.1578
$ws.Range("J91").Value = 1425.5
$ws.Range("L91").Value = 1425.5
$ws.Range("N91").Value = -4233.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 695.5833
$ws.Range("I97").Value = 658
$ws.Range("J97").Value = 808.3333
$ws.Range("K97").Value = 658
$ws.Range("L97").Value = 808.3333
$ws.Range("M97").Value = -162
$ws.Range("N97").Value = -1800.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2057.9
$ws.Range("I110").Value = 1682.7142
$ws.Range("K110").Value = 1682.7142
$ws.Range("M110").Value = 362.2858000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3896.122
$ws.Range("I132").Value = 3786.7534
$ws.Range("K132").Value = 11360.2602
$ws.Range("M132").Value = -8830.260200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6193.974
$ws.Range("I136").Value = 4043.0312
$ws.Range("K136").Value = 12129.0936
$ws.Range("M136").Value = -9579.0936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 24024.154
$ws.Range("I16").Value = 20483.625
$ws.Range("K16").Value = 20483.625
$ws.Range("M16").Value = -20196.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5742.522
$ws.Range("I31").Value = 6763.6206
$ws.Range("J31").Value = 4000.647
$ws.Range("K31").Value = 6763.6206
$ws.Range("L31").Value = 4000.647
$ws.Range("M31").Value = -6468.6206
$ws.Range("N31").Value = -4590.647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5742.522
$ws.Range("I34").Value = 6763.6206
$ws.Range("J34").Value = 4000.647
$ws.Range("K34").Value = 6763.6206
$ws.Range("L34").Value = 4000.647
$ws.Range("M34").Value = -6561.6206
$ws.Range("N34").Value = -4404.647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2016.5
$ws.Range("I58").Value = 2016.5
$ws.Range("K58").Value = 2016.5
$ws.Range("M58").Value = -1813.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3609.25
$ws.Range("I62").Value = 3625
$ws.Range("K62").Value = 3625
$ws.Range("M62").Value = -3001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3609.25
$ws.Range("I65").Value = 3625
$ws.Range("K65").Value = 18125
$ws.Range("M65").Value = -15005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3767.3333
$ws.Range("I105").Value = 3767.3333
$ws.Range("K105").Value = 3767.3333
$ws.Range("M105").Value = -2020.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 24024.154
$ws.Range("I113").Value = 20483.625
$ws.Range("K113").Value = 20483.625
$ws.Range("M113").Value = -18313.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 275520.53
$ws.Range("I122").Value = 378229.12
$ws.Range("J122").Value = 1631
$ws.Range("K122").Value = 1134687.36
$ws.Range("L122").Value = 4893
$ws.Range("M122").Value = -1132237.36
$ws.Range("N122").Value = -9793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1905.619
$ws.Range("I132").Value = 1905.4
$ws.Range("K132").Value = 5716.200000000001
$ws.Range("M132").Value = -3186.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2016.5
$ws.Range("I136").Value = 2016.5
$ws.Range("K136").Value = 6049.5
$ws.Range("M136").Value = -3499.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 170.9
$ws.Range("I12").Value = 357
$ws.Range("J12").Value = 91.14286
$ws.Range("K12").Value = 1071
$ws.Range("L12").Value = 273.42858
$ws.Range("M12").Value = -898
$ws.Range("N12").Value = -619.42858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 20000
$ws.Range("I120").Value = 20000
$ws.Range("K120").Value = 60000
$ws.Range("M120").Value = -55162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 7501.3335
$ws.Range("I22").Value = 8166.6665
$ws.Range("K22").Value = 8166.6665
$ws.Range("M22").Value = -7637.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3672.7144
$ws.Range("J80").Value = 4947.143
$ws.Range("L80").Value = 4947.143
$ws.Range("N80").Value = -6943.143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3672.7144
$ws.Range("J83").Value = 4947.143
$ws.Range("L83").Value = 24735.715
$ws.Range("N83").Value = -34719.715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4913.6587
$ws.Range("I132").Value = 4732.528
$ws.Range("J132").Value = 6217.8
$ws.Range("K132").Value = 14197.584
$ws.Range("L132").Value = 18653.4
$ws.Range("M132").Value = -11667.584
$ws.Range("N132").Value = -23713.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4109.3335
$ws.Range("J46").Value = 3757.3171
$ws.Range("L46").Value = 3757.3171
$ws.Range("N46").Value = -4133.3171

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 791.7083
$ws.Range("I55").Value = 866.46155
$ws.Range("K55").Value = 866.46155
$ws.Range("M55").Value = -693.46155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1414.4615
$ws.Range("J82").Value = 1123.75
$ws.Range("L82").Value = 1123.75
$ws.Range("N82").Value = -1845.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1414.4615
$ws.Range("J85").Value = 1123.75
$ws.Range("L85").Value = 1123.75
$ws.Range("N85").Value = -3619.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 19676
$ws.Range("I122").Value = 19676
$ws.Range("K122").Value = 59028
$ws.Range("M122").Value = -56578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2469
$ws.Range("I132").Value = 1974.0667
$ws.Range("K132").Value = 5922.2001
$ws.Range("M132").Value = -3392.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2877.7
$ws.Range("I126").Value = 2877.7
$ws.Range("K126").Value = 8633.099999999999
$ws.Range("M126").Value = -6163.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2983.5625
$ws.Range("I132").Value = 3087.5898
$ws.Range("J132").Value = 2532.7778
$ws.Range("K132").Value = 9262.769400000001
$ws.Range("L132").Value = 7598.3334
$ws.Range("M132").Value = -6732.769400000001
$ws.Range("N132").Value = -12658.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2906.5715
$ws.Range("I136").Value = 2529.64
$ws.Range("J136").Value = 3848.9
$ws.Range("K136").Value = 7588.92
$ws.Range("L136").Value = 11546.7
$ws.Range("M136").Value = -5038.92
$ws.Range("N136").Value = -16646.7
